$d = $word.ActiveDocument

function Find-ParagraphByText([string]$searchText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        return $null
    }
    return $rng.Paragraphs(1)
}

function Remove-ParagraphByText([string]$searchText) {
    $para = Find-ParagraphByText $searchText
    if ($null -eq $para) {
        return
    }
    $full = $d.Range($para.Range.Start, $para.Range.End)
    $full.Delete()
}

# ---------------------------------------------------------------------------
# 1) New intro paragraphs after "small bowel (jejunum)" (Anatomy section)
# ---------------------------------------------------------------------------
$anchor = Find-ParagraphByText("small bowel (jejunum)")
$anchorRange = $anchor.Range
$anchorRange.Collapse(0)

$newTexts = @(
    "We’ll start with reviewing some anatomy about how the body digests food.",
    "Food moves from the throat to the esophagus, and from there to the stomach.",
    "From the stomach, food moved through a valve called the pylorus into the small intestines"
)

$prevPara = $anchor
foreach ($t in $newTexts) {
    $prevPara.Range.InsertParagraphAfter() | Out-Null
    $newPara = $prevPara.Next()
    $newPara.Style = "BodyText"
    $newPara.Range.Text = $t
    $prevPara = $newPara
}

# ---------------------------------------------------------------------------
# 2) Typo fixes (simple text replace, in place)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "In many ways, these to different types of esophageal cancer behave the same.",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "In many ways, these two different types of esophageal cancer behave in similar fashion.",
    2) | Out-Null

$d.Content.Find.Execute(
    "If we look at the wasll of the esophagus, we see several layers:",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "If we look at the walls of the esophagus, we see several layers:",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Chemotherapy section: add a "Central Venous Port" slideshow link,
#    drop the old Peripheral IVs / PICC Lines / Central Venous Port (x2)
#    sub-sections, and renumber the remaining headings.
# ---------------------------------------------------------------------------

# 3a) New hyperlink paragraph right after the "Central Venous Port" bullet
$cvpBullet = Find-ParagraphByText("Central Venous Port")
$cvpRange = $cvpBullet.Range
$cvpRange.Collapse(0)
$cvpRange.InsertParagraphAfter() | Out-Null
$hlPara = $cvpBullet.Next()
$hlPara.Range.ListFormat.RemoveNumbers()
$hlPara.Style = "FirstParagraph"
$d.Hyperlinks.Add($hlPara.Range, "lci_central-venous-port.htm", "", "", "Central Venous Port") | Out-Null

# 3b) Remove the now-superseded sub-sections entirely (heading + body)
Remove-ParagraphByText("33 Peripheral IVs")
Remove-ParagraphByText("Some patients can be treated with an intravenous line placed in the hand or arm for each dose of chemotherapy. The catheter is placed at the beginning of each dose and removed that day.")

Remove-ParagraphByText("34 PICC Lines")
Remove-ParagraphByText("A PICC line is placed in Radiology and stays in place during the treatment course")

Remove-ParagraphByText("35 Central Venous Port")
Remove-ParagraphByText("A central venous port is an implantable device that makes the administration of chemotherapy easier")

Remove-ParagraphByText("36 Central Venous Port")
Remove-ParagraphByText("A central venous port is typically placed underneath the skin below the right collarbone")

Remove-ParagraphByText("37 Central Venous Port")
Remove-ParagraphByText("When it is time for chemotherapy, a needle is inserted through the skin into the port")

# 3c) Renumber the remaining headings (38/39/40 -> 33/34/35)
$d.Content.Find.Execute("38 Restaging", $true, $true, $false, $false, $false, $true, 1, $false, "33 Restaging", 2) | Out-Null
$d.Content.Find.Execute("39 Nutrition", $true, $true, $false, $false, $false, $true, 1, $false, "34 Nutrition", 2) | Out-Null
$d.Content.Find.Execute("40 Surgery", $true, $true, $false, $false, $false, $true, 1, $false, "35 Surgery", 2) | Out-Null

Write-Output "done"
